$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.361.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.794.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.005'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4510'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3599'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07085'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8835'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07735'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.758.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.281'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.327'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.007'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008529'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '26.382.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.971'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.996.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.969'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.022'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '111.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.847'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08687'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.076'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.760'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.442'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7220'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.103'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.004'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.067'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01930'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05102'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.859'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5052'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.825'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1521'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.010'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4629'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.852'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.570'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.49%  '
